# Update the confusion matrix with the results produced using the
# mse_vec output from getPerceptronWeights (for plotting of MSE).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 10,10
$arr[0,0] = 955
$arr[0,1] = 1
$arr[0,2] = 3
$arr[0,3] = 1
$arr[0,4] = 1
$arr[0,5] = 6
$arr[0,6] = 10
$arr[0,7] = 1
$arr[0,8] = 2
$arr[0,9] = 0
$arr[1,0] = 0
$arr[1,1] = 1129
$arr[1,2] = 2
$arr[1,3] = 1
$arr[1,4] = 0
$arr[1,5] = 0
$arr[1,6] = 2
$arr[1,7] = 0
$arr[1,8] = 1
$arr[1,9] = 0
$arr[2,0] = 12
$arr[2,1] = 15
$arr[2,2] = 947
$arr[2,3] = 10
$arr[2,4] = 4
$arr[2,5] = 1
$arr[2,6] = 4
$arr[2,7] = 13
$arr[2,8] = 26
$arr[2,9] = 0
$arr[3,0] = 0
$arr[3,1] = 5
$arr[3,2] = 8
$arr[3,3] = 946
$arr[3,4] = 1
$arr[3,5] = 22
$arr[3,6] = 0
$arr[3,7] = 10
$arr[3,8] = 14
$arr[3,9] = 4
$arr[4,0] = 0
$arr[4,1] = 13
$arr[4,2] = 3
$arr[4,3] = 0
$arr[4,4] = 900
$arr[4,5] = 0
$arr[4,6] = 12
$arr[4,7] = 2
$arr[4,8] = 3
$arr[4,9] = 49
$arr[5,0] = 4
$arr[5,1] = 4
$arr[5,2] = 2
$arr[5,3] = 27
$arr[5,4] = 7
$arr[5,5] = 828
$arr[5,6] = 7
$arr[5,7] = 2
$arr[5,8] = 7
$arr[5,9] = 4
$arr[6,0] = 8
$arr[6,1] = 4
$arr[6,2] = 4
$arr[6,3] = 0
$arr[6,4] = 11
$arr[6,5] = 10
$arr[6,6] = 918
$arr[6,7] = 0
$arr[6,8] = 3
$arr[6,9] = 0
$arr[7,0] = 0
$arr[7,1] = 38
$arr[7,2] = 13
$arr[7,3] = 1
$arr[7,4] = 14
$arr[7,5] = 0
$arr[7,6] = 0
$arr[7,7] = 925
$arr[7,8] = 3
$arr[7,9] = 34
$arr[8,0] = 5
$arr[8,1] = 2
$arr[8,2] = 4
$arr[8,3] = 23
$arr[8,4] = 10
$arr[8,5] = 29
$arr[8,6] = 4
$arr[8,7] = 7
$arr[8,8] = 885
$arr[8,9] = 5
$arr[9,0] = 6
$arr[9,1] = 11
$arr[9,2] = 4
$arr[9,3] = 8
$arr[9,4] = 27
$arr[9,5] = 5
$arr[9,6] = 1
$arr[9,7] = 26
$arr[9,8] = 8
$arr[9,9] = 913
$ws.Range("C3:L12").Value2 = $arr

# Move the active selection to match the state the workbook was last saved in
$ws.Range("L16").Select()

# Flag hit-ratios that clear the new 80% acceptance threshold in green
$hitRatioRange = $ws.Range("O3:O12")
$fcAbs = $hitRatioRange.FormatConditions.Add(1, 5, "=80")
$fcAbs.Font.Color = 5287936
$fcPct = $hitRatioRange.FormatConditions.Add(1, 5, "=0.8")
$fcPct.Font.Color = 5287936

Write-Host "Confusion matrix refreshed with mse_vec-derived results"
